$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of data below the existing table (rows 9 and 10)
$ws.Range("A9").Value = "Price option page check for open mandatory field"
$ws.Range("A10").Value = "Price option page check for filled mandatory field"

$ws.Range("B9").Value = "<MissingMandatoryField>"
$ws.Range("B10").Value = "<FilledMandatoryField>"

# Update B1: id=priceTable -> XPath style locator
$ws.Range("B1").Value = "//*[@id='priceTable']"

# Apply the same number format used elsewhere (e.g. A1/A2 header cells) to the new B9/B10 cells
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"

# Adjust column widths for A and B to fit the new, longer content
$ws.Columns.Item(1).ColumnWidth = 40.833333333333336
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668

# Reposition the screenshot picture further down the sheet so it still starts
# right below the data table (which has grown by two rows)
$shp = $ws.Shapes.Item(1)
$shp.Top = 197.63511811023622
$shp.Width = 865.8

# Update selection to reflect the new active cell
$ws.Range("C10").Select()
